$wb = $excel.ActiveWorkbook

# --- Step 1: rename the existing "总计" sheet to "2022-Q1" and rebuild its content ---
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"
$q1.Cells.Clear()

# Header row (B1:H1), bold/border/center style copied from an existing header cell
$styleSrcHeader = $wb.Worksheets.Item("2021-Q4").Range("B1")
$styleSrcHeader.Copy($q1.Range("B1:H1"))
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Column A (row index) style, copied from an existing data cell in column A
$styleSrcA = $wb.Worksheets.Item("2021-Q4").Range("A2")
$styleSrcA.Copy($q1.Range("A2:A38"))

# Columns B:G hold numeric-looking values stored as TEXT (matches source data)
$q1.Range("B2:G38").NumberFormat = "@"

$q1Data = @(
  @("161725","招商中证白酒指数","688.84","94.34","2.30","15.8433",10),
  @("003378","泰康策略优选灵活配置混合","19.86","81.29","5.44","1.0804",2),
  @("010874","泰康品质生活混合A","13.17","81.43","5.75","0.7573",2),
  @("002621","中欧消费主题股票A","19.29","88.29","3.41","0.6578",10),
  @("012493","长信内需均衡混合型证券投资基金A","11.45","87.51","3.80","0.4351",9),
  @("006926","长城量化精选股票A","4.97","90.86","7.56","0.3757",9),
  @("519125","浦银安盛消费升级混合A","4.64","89.16","7.03","0.3262",3),
  @("010875","泰康品质生活混合C","4.39","81.43","5.75","0.2524",2),
  @("005014","泰康景泰回报混合A","11.64","27.21","2.09","0.2433",3),
  @("009875","天弘甄选食品饮料股票A","3.67","82.33","6.09","0.2235",5),
  @("002697","中欧消费主题股票C","6.20","88.29","3.41","0.2114",10),
  @("519170","浦银安盛增长动力灵活配置混合","8.12","85.61","2.53","0.2054",7),
  @("501038","银华明择多策略定期开放混合","3.08","83.04","6.52","0.2008",6),
  @("519176","浦银安盛消费升级混合C","2.33","89.16","7.03","0.1638",3),
  @("009876","天弘甄选食品饮料股票C","2.62","82.33","6.09","0.1596",5),
  @("003190","创金合信消费主题股票A","2.31","86.40","5.30","0.1224",7),
  @("001030","天弘云端生活优选灵活配置混合","1.61","79.35","6.72","0.1082",2),
  @("003191","创金合信消费主题股票C","1.89","86.40","5.30","0.1002",7),
  @("519115","浦银安盛红利精选混合","1.79","88.68","4.90","0.0877",5),
  @("011463","长城量化精选股票C","1.09","90.86","7.56","0.0824",9),
  @("510630","华夏上证主要消费ETF","3.36","99.52","2.32","0.0780",9),
  @("233008","大摩消费领航混合基金","0.89","79.72","5.72","0.0509",8),
  @("519120","浦银安盛新兴产业混合","2.21","90.11","2.26","0.0499",10),
  @("009954","北信瑞丰优选成长股票","0.57","94.37","5.03","0.0287",6),
  @("010157","汇安中证500指数增强A","1.61","94.76","1.68","0.0270",10),
  @("012494","长信内需均衡混合型证券投资基金C","0.70","87.51","3.80","0.0266",9),
  @("010331","天弘消费股票A","0.47","83.48","5.58","0.0262",3),
  @("002159","东吴国企改革主题灵活配置混合","0.24","90.26","7.19","0.0173",9),
  @("165531","信诚多策略灵活配置混合（LOF）","0.75","68.43","2.19","0.0164",10),
  @("519172","浦银安盛睿智精选灵活配置混合A","0.49","90.86","3.25","0.0159",7),
  @("002512","长城久润混合","0.34","88.14","4.51","0.0153",8),
  @("005015","泰康景泰回报混合C","0.63","27.21","2.09","0.0132",3),
  @("010332","天弘消费股票C","0.20","83.48","5.58","0.0112",3),
  @("010158","汇安中证500指数增强C","0.64","94.76","1.68","0.0108",10),
  @("004805","长信消费精选行业量化股票","0.11","89.89","6.75","0.0074",10),
  @("519173","浦银安盛睿智精选灵活配置混合C","0.18","90.86","3.25","0.0058",7),
  @("009027","浦银安盛安远回报一年持有期混合A","0.79","20.03","0.68","0.0054",4)
)

$r = 2
foreach ($row in $q1Data) {
  $q1.Cells.Item($r, 1).Value = $r - 2
  $q1.Cells.Item($r, 2).Value = $row[0]
  $q1.Cells.Item($r, 3).Value = $row[1]
  $q1.Cells.Item($r, 4).Value = $row[2]
  $q1.Cells.Item($r, 5).Value = $row[3]
  $q1.Cells.Item($r, 6).Value = $row[4]
  $q1.Cells.Item($r, 7).Value = $row[5]
  $q1.Cells.Item($r, 8).Value = $row[6]
  $r = $r + 1
}

# --- Step 2: add a brand-new "总计" sheet after "2022-Q1" with the refreshed totals table ---
$sheets = $wb.Worksheets
$total = $sheets.Add($null, $q1)
$total.Name = "总计"
$total.Outline.SummaryBelow = $true
$total.Outline.SummaryRight = $true

$styleSrcHeader.Copy($total.Range("B1:D1"))
$styleSrcA.Copy($total.Range("A2:A7"))
$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$totalData = @(
  @("2022-Q1",37,22.04),
  @("2021-Q4",44,27.95),
  @("2021-Q3",9,1.54),
  @("2021-Q2",32,28.8),
  @("2021-Q1",39,30.32),
  @("2020-Q4",30,22.27)
)

$r = 2
foreach ($row in $totalData) {
  $total.Cells.Item($r, 1).Value = $r - 2
  $total.Cells.Item($r, 2).Value = $row[0]
  $total.Cells.Item($r, 3).Value = $row[1]
  $total.Cells.Item($r, 4).Value = $row[2]
  $r = $r + 1
}
